# Concise marksheet: update the corrected (right) answer count and the
# resulting total, and the "correct/total" summary text on the Total row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): number right went from 3 -> 5
$ws.Range("B11").Value = 5

# Row 12 ("Total"): total score went from 54 -> 90, and the
# "correct/total" label updated to match (53/84 -> 90/140)
$ws.Range("B12").Value = 90
$ws.Range("E12").Value = "90/140"
